$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.114.04"
$ws.Range("E2").Value = "  +1.50%  "

$ws.Range("D3").Value = "2.060.49"
$ws.Range("E3").Value = "  -2.19%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.662"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.85"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +16.94%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.380"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0797"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.87%  "

$ws.Range("E12").Value = "  +5.83%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.42%  "

$ws.Range("E14").Value = "  -2.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.817"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.03%  "

$ws.Range("D17").Value = "2.055.19"
$ws.Range("E17").Value = "  -2.50%  "

$ws.Range("D18").Value = "37.065.32"
$ws.Range("E18").Value = "  +1.45%  "

$ws.Range("D19").Value = "0.0₃0953"
$ws.Range("E19").Value = "  +14.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.70%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.17%  "

$ws.Range("E22").Value = "  +3.79%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.27%  "

$ws.Range("E24").Value = "  -0.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.78%  "

$ws.Range("E27").Value = "  -1.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.46%  "

$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("E30").Value = "  +0.28%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.53%  "

$ws.Range("E32").Value = "  +10.91%  "

$ws.Range("E33").Value = "  +2.42%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.25%  "

$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.73%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0853"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.61%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.66%  "

$ws.Range("E39").Value = "  +1.10%  "

$ws.Range("E40").Value = "  +26.48%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0224"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.10%  "

$ws.Range("E43").Value = "  -3.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.43%  "

$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +53.04%  "

$ws.Range("B46").Value = "HuobiToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.32%  "

$ws.Range("E47").Value = "  +7.49%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.298.17"
$ws.Range("E48").Value = "  -2.84%  "

$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.83%  "

$ws.Range("B50").Value = "Gas"
$ws.Range("C50").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "12.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -54.31%  "

$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.76%  "
